$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a numeric-looking string (e.g. "7.00", "0.102").
# Excel would otherwise auto-convert these to real numbers on assignment
# (dropping trailing zeros / switching to scientific notation), but the
# source data must stay text, exactly as printed. Mark them as Text first,
# then restore the Normal style so no stray number format lingers on the cell.
$textCells = @("D5", "D7", "D10", "D11", "D13", "D16", "D19", "D21", "D22", "D23", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D43", "D47", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '41.380.14'
$ws.Range('E2').Value = '  -3.20%  '
$ws.Range('D3').Value = '2.478.12'
$ws.Range('E3').Value = '  -2.51%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '312.29'
$ws.Range('E5').Value = '  +0.28%  '
$ws.Range('E6').Value = '  -5.77%  '
$ws.Range('D7').Value = '0.548'
$ws.Range('E7').Value = '  -3.39%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -4.38%  '
$ws.Range('D10').Value = '33.58'
$ws.Range('E10').Value = '  -5.67%  '
$ws.Range('D11').Value = '0.0783'
$ws.Range('E11').Value = '  -2.75%  '
$ws.Range('E12').Value = '  -0.74%  '
$ws.Range('D13').Value = '7.00'
$ws.Range('E13').Value = '  -4.39%  '
$ws.Range('D14').Value = '2.859.01'
$ws.Range('E14').Value = '  -2.33%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.501.69'
$ws.Range('E15').Value = '  -5.77%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = '15.19'
$ws.Range('E16').Value = '  -1.24%  '
$ws.Range('E17').Value = '  -3.21%  '
$ws.Range('D18').Value = '41.400.08'
$ws.Range('E18').Value = '  -3.12%  '
$ws.Range('D19').Value = '6.34'
$ws.Range('E19').Value = '  -5.75%  '
$ws.Range('D20').Value = '0.0₃0928'
$ws.Range('E20').Value = '  -2.72%  '
$ws.Range('D21').Value = '11.26'
$ws.Range('E21').Value = '  -8.84%  '
$ws.Range('D22').Value = '68.74'
$ws.Range('E22').Value = '  -1.95%  '
$ws.Range('D23').Value = '237.70'
$ws.Range('E23').Value = '  -2.39%  '
$ws.Range('E24').Value = '  -4.41%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').Value = '1.91'
$ws.Range('E25').Value = '  -5.94%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').Value = '24.11'
$ws.Range('E27').Value = '  -6.05%  '
$ws.Range('E28').Value = '  -4.37%  '
$ws.Range('D29').Value = '9.70'
$ws.Range('E29').Value = '  -4.69%  '
$ws.Range('D30').Value = '36.69'
$ws.Range('E30').Value = '  -4.94%  '
$ws.Range('D31').Value = '152.01'
$ws.Range('E31').Value = '  -3.77%  '
$ws.Range('D32').Value = '5.49'
$ws.Range('E32').Value = '  -6.61%  '
$ws.Range('E33').Value = '  -3.12%  '
$ws.Range('D34').Value = '2.57'
$ws.Range('E34').Value = '  -6.56%  '
$ws.Range('D35').Value = '0.0748'
$ws.Range('E35').Value = '  -5.84%  '
$ws.Range('D36').Value = '3.09'
$ws.Range('E36').Value = '  -2.33%  '
$ws.Range('D37').Value = '17.62'
$ws.Range('E37').Value = '  -2.61%  '
$ws.Range('D38').Value = '1.88'
$ws.Range('E38').Value = '  -4.68%  '
$ws.Range('E39').Value = '  -2.79%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '4.28'
$ws.Range('E40').Value = '  +3.17%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '0.102'
$ws.Range('E41').Value = '  -8.57%  '
$ws.Range('E42').Value = '  +0.22%  '
$ws.Range('D43').Value = '19.65'
$ws.Range('E43').Value = '  -10.37%  '
$ws.Range('D44').Value = '1.984.56'
$ws.Range('E44').Value = '  -0.58%  '
$ws.Range('E45').Value = '  -4.19%  '
$ws.Range('E46').Value = '  -8.20%  '
$ws.Range('D47').Value = '8.74'
$ws.Range('E47').Value = '  -3.82%  '
$ws.Range('D48').Value = '2.720.21'
$ws.Range('E48').Value = '  -2.16%  '
$ws.Range('D49').Value = '69.80'
$ws.Range('E49').Value = '  -3.75%  '
$ws.Range('D50').Value = '97.34'
$ws.Range('E50').Value = '  -4.14%  '
$ws.Range('D51').Value = '74.79'
$ws.Range('E51').Value = '  -6.70%  '

# Drop the temporary Text number format back to Normal now that the literal
# text has been committed, so cell styling matches the untouched cells.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
